# Update heating technology availability table:
# - id_heating_technology 25 / id_heating_system_action 1 (row 7): all years (2010-2050, cols E:AS) -> 0
# - id_heating_technology 25 / id_heating_system_action 2 (row 32): future years only (2020-2050, cols O:AS) -> 0
# - id_heating_technology 25 / id_heating_system_action 3 (row 57): all years (2010-2050, cols E:AS) -> 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: columns E (5) through AS (45) -> 0
$ws.Range($ws.Cells.Item(7, 5), $ws.Cells.Item(7, 45)).Value = 0

# Row 32: columns O (15) through AS (45) -> 0
$ws.Range($ws.Cells.Item(32, 15), $ws.Cells.Item(32, 45)).Value = 0

# Row 57: columns E (5) through AS (45) -> 0
$ws.Range($ws.Cells.Item(57, 5), $ws.Cells.Item(57, 45)).Value = 0

# Update the active selection to match the last edited range (O32:AS32)
$ws.Range("O32:AS32").Select()
